$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Hspg2"
$ws.Cells.Item(2, 3).Value = "Col13a1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 255.961578
$ws.Cells.Item(2, 8).Value = 511.9231559999999
$ws.Cells.Item(2, 9).Value = 0.4306976391465507
$ws.Cells.Item(2, 10).Value = 0.3553431329241539
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.4274105
$ws.Cells.Item(2, 14).Value = 0.854821
$ws.Cells.Item(2, 15).Value = 0.4422185568930271
$ws.Cells.Item(2, 16).Value = 0.3984644485961792
$ws.Cells.Item(2, 17).Value = 109.400666033769
$ws.Cells.Item(2, 18).Value = 437.602664135076
$ws.Cells.Item(2, 19).Value = 0.1904624884406214
$ws.Cells.Item(2, 20).Value = 0.1415916055230618

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Hspg2"
$ws.Cells.Item(3, 3).Value = "Col13a1"
$ws.Cells.Item(3, 4).Value = "Neutro"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 255.961578
$ws.Cells.Item(3, 8).Value = 511.9231559999999
$ws.Cells.Item(3, 9).Value = 0.4306976391465507
$ws.Cells.Item(3, 10).Value = 0.3553431329241539
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.2122596666666667
$ws.Cells.Item(3, 14).Value = 0.636779
$ws.Cells.Item(3, 15).Value = 0.2196136114576696
$ws.Cells.Item(3, 16).Value = 0.2968268129966699
$ws.Cells.Item(3, 17).Value = 54.33031922575399
$ws.Cells.Item(3, 18).Value = 325.9819153545239
$ws.Cells.Item(3, 19).Value = 0.09458706397926618
$ws.Cells.Item(3, 20).Value = 0.1054753696661287

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Hspg2"
$ws.Cells.Item(4, 3).Value = "Col13a1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 255.961578
$ws.Cells.Item(4, 8).Value = 511.9231559999999
$ws.Cells.Item(4, 9).Value = 0.4306976391465507
$ws.Cells.Item(4, 10).Value = 0.3553431329241539
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.326844
$ws.Cells.Item(4, 14).Value = 0.653688
$ws.Cells.Item(4, 15).Value = 0.3381678316493033
$ws.Cells.Item(4, 16).Value = 0.304708738407151
$ws.Cells.Item(4, 17).Value = 83.65950599983199
$ws.Cells.Item(4, 18).Value = 334.638023999328
$ws.Cells.Item(4, 19).Value = 0.1456480867266631
$ws.Cells.Item(4, 20).Value = 0.1082761577349635

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Hspg2"
$ws.Cells.Item(5, 3).Value = "Col13a1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 89.56505566666668
$ws.Cells.Item(5, 8).Value = 268.695167
$ws.Cells.Item(5, 9).Value = 0.1507080020645237
$ws.Cells.Item(5, 10).Value = 0.1865103801699463
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.4274105
$ws.Cells.Item(5, 14).Value = 0.854821
$ws.Cells.Item(5, 15).Value = 0.4422185568930271
$ws.Cells.Item(5, 16).Value = 0.3984644485961792
$ws.Cells.Item(5, 17).Value = 38.28104522501784
$ws.Cells.Item(5, 18).Value = 229.686271350107
$ws.Cells.Item(5, 19).Value = 0.06664587518520501
$ws.Cells.Item(5, 20).Value = 0.07431775579188142

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Hspg2"
$ws.Cells.Item(6, 3).Value = "Col13a1"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 89.56505566666668
$ws.Cells.Item(6, 8).Value = 268.695167
$ws.Cells.Item(6, 9).Value = 0.1507080020645237
$ws.Cells.Item(6, 10).Value = 0.1865103801699463
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.2122596666666667
$ws.Cells.Item(6, 14).Value = 0.636779
$ws.Cells.Item(6, 15).Value = 0.2196136114576696
$ws.Cells.Item(6, 16).Value = 0.2968268129966699
$ws.Cells.Item(6, 17).Value = 19.01104886078811
$ws.Cells.Item(6, 18).Value = 171.099439747093
$ws.Cells.Item(6, 19).Value = 0.03309752860895997
$ws.Cells.Item(6, 20).Value = 0.05536128173664247

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Hspg2"
$ws.Cells.Item(7, 3).Value = "Col13a1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 89.56505566666668
$ws.Cells.Item(7, 8).Value = 268.695167
$ws.Cells.Item(7, 9).Value = 0.1507080020645237
$ws.Cells.Item(7, 10).Value = 0.1865103801699463
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.326844
$ws.Cells.Item(7, 14).Value = 0.653688
$ws.Cells.Item(7, 15).Value = 0.3381678316493033
$ws.Cells.Item(7, 16).Value = 0.304708738407151
$ws.Cells.Item(7, 17).Value = 29.27380105431601
$ws.Cells.Item(7, 18).Value = 175.642806325896
$ws.Cells.Item(7, 19).Value = 0.05096459827035869
$ws.Cells.Item(7, 20).Value = 0.05683134264142246

$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Hspg2"
$ws.Cells.Item(8, 3).Value = "Col13a1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.1344533333333333
$ws.Cells.Item(8, 8).Value = 0.4033600000000001
$ws.Cells.Item(8, 9).Value = 0.0002262399446609557
$ws.Cells.Item(8, 10).Value = 0.0002799857838356637
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.4274105
$ws.Cells.Item(8, 14).Value = 0.854821
$ws.Cells.Item(8, 15).Value = 0.4422185568930271
$ws.Cells.Item(8, 16).Value = 0.3984644485961792
$ws.Cells.Item(8, 17).Value = 0.05746676642666668
$ws.Cells.Item(8, 18).Value = 0.3448005985600001
$ws.Cells.Item(8, 19).Value = 0.0001000475018395262
$ws.Cells.Item(8, 20).Value = 0.0001115643809708468

$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Hspg2"
$ws.Cells.Item(9, 3).Value = "Col13a1"
$ws.Cells.Item(9, 4).Value = "Neutro"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.1344533333333333
$ws.Cells.Item(9, 8).Value = 0.4033600000000001
$ws.Cells.Item(9, 9).Value = 0.0002262399446609557
$ws.Cells.Item(9, 10).Value = 0.0002799857838356637
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.2122596666666667
$ws.Cells.Item(9, 14).Value = 0.636779
$ws.Cells.Item(9, 15).Value = 0.2196136114576696
$ws.Cells.Item(9, 16).Value = 0.2968268129966699
$ws.Cells.Item(9, 17).Value = 0.02853901971555555
$ws.Cells.Item(9, 18).Value = 0.25685117744
$ws.Cells.Item(9, 19).Value = 0.00004968537130297581
$ws.Cells.Item(9, 20).Value = 0.0000831072879003146

$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Hspg2"
$ws.Cells.Item(10, 3).Value = "Col13a1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.1344533333333333
$ws.Cells.Item(10, 8).Value = 0.4033600000000001
$ws.Cells.Item(10, 9).Value = 0.0002262399446609557
$ws.Cells.Item(10, 10).Value = 0.0002799857838356637
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.326844
$ws.Cells.Item(10, 14).Value = 0.653688
$ws.Cells.Item(10, 15).Value = 0.3381678316493033
$ws.Cells.Item(10, 16).Value = 0.304708738407151
$ws.Cells.Item(10, 17).Value = 0.04394526528000001
$ws.Cells.Item(10, 18).Value = 0.26367159168
$ws.Cells.Item(10, 19).Value = 0.00007650707151845377
$ws.Cells.Item(10, 20).Value = 0.00008531411496450237

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Hspg2"
$ws.Cells.Item(11, 3).Value = "Col13a1"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.285862
$ws.Cells.Item(11, 8).Value = 0.8575860000000001
$ws.Cells.Item(11, 9).Value = 0.0004810100386305295
$ws.Cells.Item(11, 10).Value = 0.0005952793742971328
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.4274105
$ws.Cells.Item(11, 14).Value = 0.854821
$ws.Cells.Item(11, 15).Value = 0.4422185568930271
$ws.Cells.Item(11, 16).Value = 0.3984644485961792
$ws.Cells.Item(11, 17).Value = 0.122180420351
$ws.Cells.Item(11, 18).Value = 0.7330825221060001
$ws.Cells.Item(11, 19).Value = 0.000212711565134252
$ws.Cells.Item(11, 20).Value = 0.0002371976676399856

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Hspg2"
$ws.Cells.Item(12, 3).Value = "Col13a1"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.285862
$ws.Cells.Item(12, 8).Value = 0.8575860000000001
$ws.Cells.Item(12, 9).Value = 0.0004810100386305295
$ws.Cells.Item(12, 10).Value = 0.0005952793742971328
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.2122596666666667
$ws.Cells.Item(12, 14).Value = 0.636779
$ws.Cells.Item(12, 15).Value = 0.2196136114576696
$ws.Cells.Item(12, 16).Value = 0.2968268129966699
$ws.Cells.Item(12, 17).Value = 0.06067697283266666
$ws.Cells.Item(12, 18).Value = 0.546092755494
$ws.Cells.Item(12, 19).Value = 0.0001056363517310437
$ws.Cells.Item(12, 20).Value = 0.0001766948795152697

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Hspg2"
$ws.Cells.Item(13, 3).Value = "Col13a1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.285862
$ws.Cells.Item(13, 8).Value = 0.8575860000000001
$ws.Cells.Item(13, 9).Value = 0.0004810100386305295
$ws.Cells.Item(13, 10).Value = 0.0005952793742971328
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.326844
$ws.Cells.Item(13, 14).Value = 0.653688
$ws.Cells.Item(13, 15).Value = 0.3381678316493033
$ws.Cells.Item(13, 16).Value = 0.304708738407151
$ws.Cells.Item(13, 17).Value = 0.09343227952800001
$ws.Cells.Item(13, 18).Value = 0.5605936771680001
$ws.Cells.Item(13, 19).Value = 0.0001626621217652338
$ws.Cells.Item(13, 20).Value = 0.0001813868271418775

$ws.Cells.Item(14, 1).Value = "Neutro"
$ws.Cells.Item(14, 2).Value = "Hspg2"
$ws.Cells.Item(14, 3).Value = "Col13a1"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 162.0686236666667
$ws.Cells.Item(14, 8).Value = 486.205871
$ws.Cells.Item(14, 9).Value = 0.2727072326181867
$ws.Cells.Item(14, 10).Value = 0.3374918978020542
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.4274105
$ws.Cells.Item(14, 14).Value = 0.854821
$ws.Cells.Item(14, 15).Value = 0.4422185568930271
$ws.Cells.Item(14, 16).Value = 0.3984644485961792
$ws.Cells.Item(14, 17).Value = 69.26983147568184
$ws.Cells.Item(14, 18).Value = 415.618988854091
$ws.Cells.Item(14, 19).Value = 0.1205961988627056
$ws.Cells.Item(14, 20).Value = 0.1344785229633736

$ws.Cells.Item(15, 1).Value = "Neutro"
$ws.Cells.Item(15, 2).Value = "Hspg2"
$ws.Cells.Item(15, 3).Value = "Col13a1"
$ws.Cells.Item(15, 4).Value = "Neutro"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 162.0686236666667
$ws.Cells.Item(15, 8).Value = 486.205871
$ws.Cells.Item(15, 9).Value = 0.2727072326181867
$ws.Cells.Item(15, 10).Value = 0.3374918978020542
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.2122596666666667
$ws.Cells.Item(15, 14).Value = 0.636779
$ws.Cells.Item(15, 15).Value = 0.2196136114576696
$ws.Cells.Item(15, 16).Value = 0.2968268129966699
$ws.Cells.Item(15, 17).Value = 34.40063203661211
$ws.Cells.Item(15, 18).Value = 309.605688329509
$ws.Cells.Item(15, 19).Value = 0.05989022022590677
$ws.Cells.Item(15, 20).Value = 0.1001766444367816

$ws.Cells.Item(16, 1).Value = "Neutro"
$ws.Cells.Item(16, 2).Value = "Hspg2"
$ws.Cells.Item(16, 3).Value = "Col13a1"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 162.0686236666667
$ws.Cells.Item(16, 8).Value = 486.205871
$ws.Cells.Item(16, 9).Value = 0.2727072326181867
$ws.Cells.Item(16, 10).Value = 0.3374918978020542
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.326844
$ws.Cells.Item(16, 14).Value = 0.653688
$ws.Cells.Item(16, 15).Value = 0.3381678316493033
$ws.Cells.Item(16, 16).Value = 0.304708738407151
$ws.Cells.Item(16, 17).Value = 52.97115723370801
$ws.Cells.Item(16, 18).Value = 317.826943402248
$ws.Cells.Item(16, 19).Value = 0.09222081352957436
$ws.Cells.Item(16, 20).Value = 0.1028367304018991

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Hspg2"
$ws.Cells.Item(17, 3).Value = "Col13a1"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 86.2797165
$ws.Cells.Item(17, 8).Value = 172.559433
$ws.Cells.Item(17, 9).Value = 0.1451798761874476
$ws.Cells.Item(17, 10).Value = 0.1197793239457127
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.4274105
$ws.Cells.Item(17, 14).Value = 0.854821
$ws.Cells.Item(17, 15).Value = 0.4422185568930271
$ws.Cells.Item(17, 16).Value = 0.3984644485961792
$ws.Cells.Item(17, 17).Value = 36.87685676912326
$ws.Cells.Item(17, 18).Value = 147.507427076493
$ws.Cells.Item(17, 19).Value = 0.06420123533752142
$ws.Cells.Item(17, 20).Value = 0.04772780226925155

$ws.Cells.Item(18, 1).Value = "sCs"
$ws.Cells.Item(18, 2).Value = "Hspg2"
$ws.Cells.Item(18, 3).Value = "Col13a1"
$ws.Cells.Item(18, 4).Value = "Neutro"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 86.2797165
$ws.Cells.Item(18, 8).Value = 172.559433
$ws.Cells.Item(18, 9).Value = 0.1451798761874476
$ws.Cells.Item(18, 10).Value = 0.1197793239457127
$ws.Cells.Item(18, 11).Value = 2
$ws.Cells.Item(18, 12).Value = 0.6666666666666666
$ws.Cells.Item(18, 13).Value = 0.2122596666666667
$ws.Cells.Item(18, 14).Value = 0.636779
$ws.Cells.Item(18, 15).Value = 0.2196136114576696
$ws.Cells.Item(18, 16).Value = 0.2968268129966699
$ws.Cells.Item(18, 17).Value = 18.3137038643845
$ws.Cells.Item(18, 18).Value = 109.882223186307
$ws.Cells.Item(18, 19).Value = 0.03188347692050269
$ws.Cells.Item(18, 20).Value = 0.03555371498970162

$ws.Cells.Item(19, 1).Value = "sCs"
$ws.Cells.Item(19, 2).Value = "Hspg2"
$ws.Cells.Item(19, 3).Value = "Col13a1"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 86.2797165
$ws.Cells.Item(19, 8).Value = 172.559433
$ws.Cells.Item(19, 9).Value = 0.1451798761874476
$ws.Cells.Item(19, 10).Value = 0.1197793239457127
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.326844
$ws.Cells.Item(19, 14).Value = 0.653688
$ws.Cells.Item(19, 15).Value = 0.3381678316493033
$ws.Cells.Item(19, 16).Value = 0.304708738407151
$ws.Cells.Item(19, 17).Value = 28.200007659726
$ws.Cells.Item(19, 18).Value = 112.800030638904
$ws.Cells.Item(19, 19).Value = 0.04909516392942347
$ws.Cells.Item(19, 20).Value = 0.03649780668675957
